$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2162.5
$ws.Range("I31").Value = 2162.5
$ws.Range("K31").Value = 6487.5
$ws.Range("M31").Value = -6257.5
$ws.Range("H42").Value = 272.5
$ws.Range("I42").Value = 38.333332
$ws.Range("J42").Value = 975
$ws.Range("K42").Value = 114.999996
$ws.Range("L42").Value = 2925
$ws.Range("M42").Value = 115.000004
$ws.Range("N42").Value = -3385
$ws.Range("H86").Value = 4304.0557
$ws.Range("I86").Value = 1694.6
$ws.Range("K86").Value = 1694.6
$ws.Range("M86").Value = -571.5999999999999
$ws.Range("H87").Value = 69650.11
$ws.Range("J87").Value = 87141.836
$ws.Range("L87").Value = 87141.836
$ws.Range("N87").Value = -89637.836
$ws.Range("H89").Value = 4304.0557
$ws.Range("I89").Value = 1694.6
$ws.Range("K89").Value = 8473
$ws.Range("M89").Value = -2857
$ws.Range("H90").Value = 69650.11
$ws.Range("J90").Value = 87141.836
$ws.Range("L90").Value = 261425.508
$ws.Range("N90").Value = -273905.508
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").Value = $null
$ws.Range("H112").Value = 1416.871
$ws.Range("I112").Value = 855
$ws.Range("J112").Value = 1455.6207
$ws.Range("K112").Value = 2565
$ws.Range("L112").Value = 4366.8621
$ws.Range("M112").Value = -1457
$ws.Range("N112").Value = -6582.8621
$ws.Range("H136").Value = 94332
$ws.Range("J136").Value = 94332
$ws.Range("L136").Value = 94332
$ws.Range("N136").Value = -104532
$ws.Range("H137").Value = 2526.72
$ws.Range("I137").Value = 1349.7778
$ws.Range("K137").Value = 4049.3334
$ws.Range("M137").Value = -1499.3334
$ws.Range("H138").Value = 2098.8936
$ws.Range("J138").Value = 3302.5625
$ws.Range("L138").Value = 9907.6875
$ws.Range("N138").Value = -20187.6875
$ws.Range("H139").Value = 95000
$ws.Range("J139").Value = 95000
$ws.Range("L139").Value = 95000
$ws.Range("N139").Value = -105280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 7224.5
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 1000
$ws.Range("N29").Value = -1616
$ws.Range("H32").Value = 3353.1406
$ws.Range("I32").Value = 3389.776
$ws.Range("K32").Value = 3389.776
$ws.Range("M32").Value = -3102.776
$ws.Range("H45").Value = 3046.56
$ws.Range("I45").Value = 4060.3572
$ws.Range("K45").Value = 4060.3572
$ws.Range("M45").Value = -3683.3572
$ws.Range("H88").Value = 1774.9
$ws.Range("I88").Value = 1458.3334
$ws.Range("J88").Value = 2249.75
$ws.Range("K88").Value = 1458.3334
$ws.Range("L88").Value = 2249.75
$ws.Range("M88").Value = -1052.3334
$ws.Range("N88").Value = -3061.75
$ws.Range("H91").Value = 1774.9
$ws.Range("I91").Value = 1458.3334
$ws.Range("J91").Value = 2249.75
$ws.Range("K91").Value = 1458.3334
$ws.Range("L91").Value = 2249.75
$ws.Range("M91").Value = -54.33339999999998
$ws.Range("N91").Value = -5057.75
$ws.Range("H122").Value = 3109.2683
$ws.Range("I122").Value = 2724.889
$ws.Range("K122").Value = 8174.667
$ws.Range("M122").Value = -5724.667
$ws.Range("H139").Value = 88571
$ws.Range("J139").Value = 88571
$ws.Range("L139").Value = 88571
$ws.Range("N139").Value = -98851

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 2000
$ws.Range("I36").Value = 2000
$ws.Range("K36").Value = 2000
$ws.Range("M36").Value = -1466
$ws.Range("H92").Value = 65396.75
$ws.Range("J92").Value = 65396.75
$ws.Range("L92").Value = 65396.75
$ws.Range("N92").Value = -70388.75
$ws.Range("H94").Value = 1076.5883
$ws.Range("I94").Value = 1153.3334
$ws.Range("J94").Value = 780.5714
$ws.Range("K94").Value = 1153.3334
$ws.Range("L94").Value = 780.5714
$ws.Range("M94").Value = -702.3334
$ws.Range("N94").Value = -1682.5714
$ws.Range("H105").Value = 1310.2258
$ws.Range("I105").Value = 1325.2916
$ws.Range("K105").Value = 1325.2916
$ws.Range("M105").Value = 421.7084
$ws.Range("H107").Value = 3078.6897
$ws.Range("I107").Value = 2854.7778
$ws.Range("K107").Value = 2854.7778
$ws.Range("M107").Value = -934.7777999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 67184.21000000001
$ws.Range("H71").Value = 67184.21000000001
$ws.Range("H92").Value = 66480.60000000001
$ws.Range("J92").Value = 66480.60000000001
$ws.Range("L92").Value = 66480.60000000001
$ws.Range("N92").Value = -71472.60000000001
$ws.Range("H134").Value = 1732.5
$ws.Range("I134").Value = 1732.5
$ws.Range("K134").Value = 5197.5
$ws.Range("M134").Value = -2662.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 436.33334
$ws.Range("I19").Value = 354.5
$ws.Range("J19").Value = 600
$ws.Range("K19").Value = 1063.5
$ws.Range("L19").Value = 1800
$ws.Range("M19").Value = -889.5
$ws.Range("N19").Value = -2148
$ws.Range("H37").Value = 84999.75
$ws.Range("J37").Value = 84999.75
$ws.Range("L37").Value = 254999.25
$ws.Range("N37").Value = -255223.25
$ws.Range("H41").Value = 255.77777
$ws.Range("J41").Value = 255.77777
$ws.Range("L41").Value = 767.33331
$ws.Range("N41").Value = -1443.33331
$ws.Range("H56").Value = 8049.375
$ws.Range("I56").Value = 8049.375
$ws.Range("K56").Value = 8049.375
$ws.Range("M56").Value = -7519.375
$ws.Range("H68").Value = 13894.286
$ws.Range("J68").Value = 15835.833
$ws.Range("L68").Value = 47507.499
$ws.Range("N68").Value = -49129.499
$ws.Range("H71").Value = 13894.286
$ws.Range("J71").Value = 15835.833
$ws.Range("L71").Value = 142522.497
$ws.Range("N71").Value = -150634.497
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null
$ws.Range("H105").Value = 13760.583
$ws.Range("J105").Value = 13760.583
$ws.Range("L105").Value = 41281.749
$ws.Range("N105").Value = -46523.749
$ws.Range("H113").Value = 1083.4286
$ws.Range("J113").Value = 952
$ws.Range("L113").Value = 2856
$ws.Range("N113").Value = -7196
$ws.Range("H136").Value = 6590.5713
$ws.Range("I136").Value = 6590.5713
$ws.Range("K136").Value = 19771.7139
$ws.Range("M136").Value = -14671.7139
$ws.Range("H138").Value = 3602.2
$ws.Range("I138").Value = 3602.2
$ws.Range("K138").Value = 10806.6
$ws.Range("M138").Value = -5666.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8072.4287
$ws.Range("I70").Value = 7599.8
$ws.Range("K70").Value = 7599.8
$ws.Range("M70").Value = -7329.8
$ws.Range("H73").Value = 8072.4287
$ws.Range("I73").Value = 7599.8
$ws.Range("K73").Value = 7599.8
$ws.Range("M73").Value = -6663.8
$ws.Range("H132").Value = 1667
$ws.Range("I132").Value = 1800.4
$ws.Range("K132").Value = 5401.200000000001
$ws.Range("M132").Value = -2871.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3086.1428
$ws.Range("I7").Value = 2800.5
$ws.Range("J7").Value = 4800
$ws.Range("K7").Value = 2800.5
$ws.Range("L7").Value = 4800
$ws.Range("M7").Value = -2688.5
$ws.Range("N7").Value = -5024
$ws.Range("H25").Value = 15555.5
$ws.Range("I25").Value = 15555.5
$ws.Range("K25").Value = 15555.5
$ws.Range("M25").Value = -15325.5
$ws.Range("H46").Value = 2667.6333
$ws.Range("I46").Value = 1100.8182
$ws.Range("J46").Value = 3574.7368
$ws.Range("K46").Value = 1100.8182
$ws.Range("L46").Value = 3574.7368
$ws.Range("M46").Value = -912.8181999999999
$ws.Range("N46").Value = -3950.7368
$ws.Range("H94").Value = 55833.332
$ws.Range("J94").Value = 55833.332
$ws.Range("L94").Value = 55833.332
$ws.Range("N94").Value = -57185.332
$ws.Range("H126").Value = 3086.1428
$ws.Range("I126").Value = 2800.5
$ws.Range("J126").Value = 4800
$ws.Range("K126").Value = 8401.5
$ws.Range("L126").Value = 14400
$ws.Range("M126").Value = -5931.5
$ws.Range("N126").Value = -19340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1036.375
$ws.Range("I81").Value = 848.5
$ws.Range("J81").Value = 1600
$ws.Range("K81").Value = 1697
$ws.Range("L81").Value = 3200
$ws.Range("M81").Value = -636
$ws.Range("N81").Value = -5322
$ws.Range("H84").Value = 1036.375
$ws.Range("I84").Value = 848.5
$ws.Range("J84").Value = 1600
$ws.Range("K84").Value = 8485
$ws.Range("L84").Value = 16000
$ws.Range("M84").Value = -3181
$ws.Range("N84").Value = -26608
$ws.Range("H126").Value = 7128.143
$ws.Range("I126").Value = 6699.6
$ws.Range("K126").Value = 20098.8
$ws.Range("M126").Value = -17628.8
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null
